$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 7787.6665
$ws.Range("J32").Value = 7439
$ws.Range("L32").Value = 7439
$ws.Range("N32").Value = -8091
$ws.Range("H74").Value = 38172.715
$ws.Range("I74").Value = 38172.715
$ws.Range("K74").Value = 38172.715
$ws.Range("M74").Value = -37236.715
$ws.Range("H76").Value = 5199.5
$ws.Range("I76").Value = 2300
$ws.Range("K76").Value = 2300
$ws.Range("M76").Value = -1985
$ws.Range("H77").Value = 38172.715
$ws.Range("I77").Value = 38172.715
$ws.Range("K77").Value = 190863.575
$ws.Range("M77").Value = -186183.575
$ws.Range("H79").Value = 5199.5
$ws.Range("I79").Value = 2300
$ws.Range("K79").Value = 2300
$ws.Range("M79").Value = -1208
$ws.Range("H92").Value = 1353.7333
$ws.Range("I92").Value = 1138.9231
$ws.Range("J92").Value = 2750
$ws.Range("K92").Value = 1138.9231
$ws.Range("L92").Value = 2750
$ws.Range("M92").Value = 109.0769
$ws.Range("N92").Value = -5246
$ws.Range("H116").Value = 21658.2
$ws.Range("I116").Value = 32781
$ws.Range("K116").Value = 32781
$ws.Range("M116").Value = -29339
$ws.Range("H138").Value = 2622.7046
$ws.Range("J138").Value = 2634.6333
$ws.Range("L138").Value = 7903.8999
$ws.Range("N138").Value = -18183.8999
$ws.Range("H141").Value = 5987
$ws.Range("I141").Value = 3793.8333
$ws.Range("J141").Value = 8618.799999999999
$ws.Range("K141").Value = 11381.4999
$ws.Range("L141").Value = 25856.4
$ws.Range("M141").Value = -6201.499899999999
$ws.Range("N141").Value = -36216.39999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 137.66667
$ws.Range("I5").Value = 146
$ws.Range("J5").Value = 127.25
$ws.Range("K5").Value = 146
$ws.Range("L5").Value = 127.25
$ws.Range("M5").Value = -34
$ws.Range("N5").Value = -351.25
$ws.Range("H32").Value = 5831.3623
$ws.Range("I32").Value = 2655.3962
$ws.Range("J32").Value = 39496.6
$ws.Range("K32").Value = 2655.3962
$ws.Range("L32").Value = 39496.6
$ws.Range("M32").Value = -2368.3962
$ws.Range("N32").Value = -40070.6
$ws.Range("H45").Value = 2287.423
$ws.Range("I45").Value = 1776.6666
$ws.Range("K45").Value = 1776.6666
$ws.Range("M45").Value = -1399.6666
$ws.Range("H61").Value = 3844.7778
$ws.Range("I61").Value = 2646.3845
$ws.Range("J61").Value = 5484.684
$ws.Range("K61").Value = 2646.3845
$ws.Range("L61").Value = 5484.684
$ws.Range("M61").Value = -2434.3845
$ws.Range("N61").Value = -5908.684
$ws.Range("H74").Value = 3860.5
$ws.Range("I74").Value = 2532.1904
$ws.Range("J74").Value = 6959.8887
$ws.Range("K74").Value = 2532.1904
$ws.Range("L74").Value = 6959.8887
$ws.Range("M74").Value = -1658.1904
$ws.Range("N74").Value = -8707.8887
$ws.Range("H77").Value = 3860.5
$ws.Range("I77").Value = 2532.1904
$ws.Range("J77").Value = 6959.8887
$ws.Range("K77").Value = 12660.952
$ws.Range("L77").Value = 34799.4435
$ws.Range("M77").Value = -8292.951999999999
$ws.Range("N77").Value = -43535.4435
$ws.Range("H122").Value = 2325.5557
$ws.Range("I122").Value = 1984.9032
$ws.Range("K122").Value = 5954.7096
$ws.Range("M122").Value = -3504.7096
$ws.Range("H132").Value = 1711.2941
$ws.Range("I132").Value = 1054.3549
$ws.Range("J132").Value = 8499.666999999999
$ws.Range("K132").Value = 3163.0647
$ws.Range("L132").Value = 25499.001
$ws.Range("M132").Value = -633.0646999999999
$ws.Range("N132").Value = -30559.001
$ws.Range("H136").Value = 3844.7778
$ws.Range("I136").Value = 2646.3845
$ws.Range("J136").Value = 5484.684
$ws.Range("K136").Value = 7939.1535
$ws.Range("L136").Value = 16454.052
$ws.Range("M136").Value = -5389.1535
$ws.Range("N136").Value = -21554.052

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 137.66667
$ws.Range("I4").Value = 146
$ws.Range("J4").Value = 127.25
$ws.Range("K4").Value = 146
$ws.Range("L4").Value = 127.25
$ws.Range("M4").Value = -31
$ws.Range("N4").Value = -357.25
$ws.Range("H22").Value = 1214.909
$ws.Range("I22").Value = 912.3333
$ws.Range("K22").Value = 912.3333
$ws.Range("M22").Value = -739.3333
$ws.Range("H94").Value = 1389.5714
$ws.Range("I94").Value = 1201
$ws.Range("J94").Value = 1861
$ws.Range("K94").Value = 1201
$ws.Range("L94").Value = 1861
$ws.Range("M94").Value = -750
$ws.Range("N94").Value = -2763
$ws.Range("H134").Value = 4324.3335
$ws.Range("I134").Value = 3223.4666
$ws.Range("J134").Value = 15333
$ws.Range("K134").Value = 9670.399800000001
$ws.Range("L134").Value = 45999
$ws.Range("M134").Value = -7135.399800000001
$ws.Range("N134").Value = -51069

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 642.35297
$ws.Range("I7").Value = 649.4167
$ws.Range("K7").Value = 649.4167
$ws.Range("M7").Value = -536.4167
$ws.Range("H57").Value = 99993.5
$ws.Range("J57").Value = 99993.5
$ws.Range("L57").Value = 99993.5
$ws.Range("N57").Value = -101113.5
$ws.Range("H132").Value = 6219.0835
$ws.Range("I132").Value = 4304.263
$ws.Range("K132").Value = 12912.789
$ws.Range("M132").Value = -10382.789

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1368.7142
$ws.Range("I5").Value = 1194
$ws.Range("K5").Value = 3582
$ws.Range("M5").Value = -3470
$ws.Range("H8").Value = 572.38464
$ws.Range("I8").Value = 572.38464
$ws.Range("K8").Value = 1717.15392
$ws.Range("M8").Value = -1578.15392
$ws.Range("H107").Value = 598
$ws.Range("J107").Value = 705.1429000000001
$ws.Range("L107").Value = 2115.4287
$ws.Range("N107").Value = -5955.4287
$ws.Range("H116").Value = 8364.333000000001
$ws.Range("I116").Value = 8400
$ws.Range("K116").Value = 25200
$ws.Range("M116").Value = -21758
$ws.Range("H135").Value = 1368.7142
$ws.Range("I135").Value = 1194
$ws.Range("K135").Value = 10746
$ws.Range("M135").Value = -8211

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 24499.5
$ws.Range("J38").Value = 24499.5
$ws.Range("L38").Value = 24499.5
$ws.Range("N38").Value = -25425.5
$ws.Range("H102").Value = 2350.45
$ws.Range("I102").Value = 1063.125
$ws.Range("K102").Value = 1063.125
$ws.Range("M102").Value = 558.875
$ws.Range("H126").Value = 7665
$ws.Range("I126").Value = 7665
$ws.Range("K126").Value = 22995
$ws.Range("M126").Value = -20525
$ws.Range("H132").Value = 5394.2593
$ws.Range("I132").Value = 3251.8
$ws.Range("J132").Value = 11515.571
$ws.Range("K132").Value = 9755.400000000001
$ws.Range("L132").Value = 34546.713
$ws.Range("M132").Value = -7225.400000000001
$ws.Range("N132").Value = -39606.713

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10041.571
$ws.Range("J7").Value = 12247.75
$ws.Range("L7").Value = 12247.75
$ws.Range("N7").Value = -12471.75
$ws.Range("H61").Value = 4246.926
$ws.Range("I61").Value = 1664.8889
$ws.Range("K61").Value = 1664.8889
$ws.Range("M61").Value = -1462.8889
$ws.Range("H113").Value = 4246.926
$ws.Range("I113").Value = 1664.8889
$ws.Range("K113").Value = 1664.8889
$ws.Range("M113").Value = 505.1111000000001
$ws.Range("H122").Value = 5185.9
$ws.Range("I122").Value = 4844.385
$ws.Range("J122").Value = 5820.143
$ws.Range("K122").Value = 14533.155
$ws.Range("L122").Value = 17460.429
$ws.Range("M122").Value = -12083.155
$ws.Range("N122").Value = -22360.429
$ws.Range("H126").Value = 10041.571
$ws.Range("J126").Value = 12247.75
$ws.Range("L126").Value = 36743.25
$ws.Range("N126").Value = -41683.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 28100
$ws.Range("J4").Value = 27400
$ws.Range("L4").Value = 27400
$ws.Range("N4").Value = -27626
$ws.Range("H113").Value = 320.07144
$ws.Range("I113").Value = 215.45833
$ws.Range("J113").Value = 947.75
$ws.Range("K113").Value = 646.37499
$ws.Range("L113").Value = 2843.25
$ws.Range("M113").Value = 1523.62501
$ws.Range("N113").Value = -7183.25
$ws.Range("H126").Value = 5834.4165
$ws.Range("I126").Value = 6046.636
$ws.Range("K126").Value = 18139.908
$ws.Range("M126").Value = -15669.908
$ws.Range("H132").Value = 2854.0356
$ws.Range("I132").Value = 2511.9443
$ws.Range("J132").Value = 3469.8
$ws.Range("K132").Value = 7535.8329
$ws.Range("L132").Value = 10409.4
$ws.Range("M132").Value = -5005.8329
$ws.Range("N132").Value = -15469.4
